$wb = $excel.ActiveWorkbook

$wsRun  = $wb.Worksheets.Item("RunManager")
$wsData = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# RunManager sheet: rename the two tests from the old Selenium login /
# homepage scenarios to the new product-search scenarios.
# ---------------------------------------------------------------------
$wsRun.Range("A2").Value = "searchProductTest1"
$wsRun.Range("B2").Value = "Search Functionality"
$wsRun.Range("A3").Value = "searchProductTest2"
$wsRun.Range("B3").Value = "Search Functionality"

$wsRun.Range("C3").Select()

# ---------------------------------------------------------------------
# Data sheet: swap username/password columns for a productName column,
# rename the tests, normalize the browser to chrome for both rows, and
# drop the (now unused) password column entirely.
# ---------------------------------------------------------------------
$wsData.Range("A2").Value = "searchProductTest1"
$wsData.Range("A3").Value = "searchProductTest2"
$wsData.Range("C3").Value = "chrome"

$wsData.Range("D1").Value = "productName"
$wsData.Range("D2").ClearFormats()
$wsData.Range("D2").Value = "'Oneplus 9R"
$wsData.Range("D3").ClearFormats()
$wsData.Range("D3").Value = "'iphone SE"

$wsData.Columns.Item(5).Delete()

$wsData.Columns.Item(1).ColumnWidth = 16.5
$wsData.Columns.Item(4).ColumnWidth = 12

$wsData.PageSetup.Orientation = 1

$wsData.Range("B3").Select()
